$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns K (actual op) and N (status)
$ws.Range("K2").Value = "actual op"
$ws.Range("N2").Value = "status"

# Row 4
$ws.Range("K4").Value = "marketing"
$ws.Range("N4").Value = "pass"

# Row 5
$ws.Range("K5").Value = "accounts"
$ws.Range("N5").Value = "pass"

# Row 6
$ws.Range("K6").Value = "sales"
$ws.Range("N6").Value = "pass"

# Row 7 - E7 changed from maths to art, H7 changed text
$ws.Range("E7").Value = "art"
$ws.Range("H7").Value = "your marks do not qualify "
$ws.Range("K7").Value = "your marks do not qualify "
$ws.Range("N7").Value = "pass"

# Row 8 - H8 changed text
$ws.Range("H8").Value = "marks do not qualify "
$ws.Range("K8").Value = "marks do not qualify "
$ws.Range("N8").Value = "pass"

# Row 9 - H9 changed text
$ws.Range("H9").Value = "marks do not qualify  "
$ws.Range("K9").Value = "marks do not qualify  "
$ws.Range("N9").Value = "pass"

# Row 10 - new row
$ws.Range("A10").Value = ">35"
$ws.Range("B10").Value = ">90"
$ws.Range("C10").Value = ">90"
$ws.Range("D10").Value = "ece"
$ws.Range("E10").Value = "maths"
$ws.Range("H10").Value = "preference do not match"
$ws.Range("K10").Value = "preference do not match"
$ws.Range("N10").Value = "pass"

# Row 11 - new row
$ws.Range("A11").Value = ">95"
$ws.Range("B11").Value = ">35"
$ws.Range("C11").Value = ">35"
$ws.Range("D11").Value = "bcom"
$ws.Range("E11").Value = "art"
$ws.Range("H11").Value = "preference do not match"
$ws.Range("K11").Value = "preference do not match"
$ws.Range("N11").Value = "pass"

# Row 12 - new row
$ws.Range("A12").Value = ">90"
$ws.Range("B12").Value = ">35"
$ws.Range("C12").Value = ">90"
$ws.Range("D12").Value = "mech"
$ws.Range("E12").Value = "art"
$ws.Range("H12").Value = "preference do not match"
$ws.Range("K12").Value = "preference do not match"
$ws.Range("N12").Value = "pass"

# Update selection to match target (O10)
$ws.Range("O10").Select()
